$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the four "Terms*" sheets to their new "TermsPage*" names.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("TermsDrugs").Name    = "TermsPageDrugs"
$wb.Worksheets.Item("TermsEnglish").Name  = "TermsPageEnglish"
$wb.Worksheets.Item("TermsGenetics").Name = "TermsPageGenetics"
$wb.Worksheets.Item("TermsSpanish").Name  = "TermsPageSpanish"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "SearchTerms" sheet right after "Definitions"
#    (i.e. it becomes the second sheet in the workbook).
# ---------------------------------------------------------------------------
$searchTerms = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Definitions"))
$searchTerms.Name = "SearchTerms"

# Header row (bold / shaded, matching the other dictionary sheets).
$searchTerms.Range("A1").Value = "Path"
$searchTerms.Range("B1").Value = "SearchTerm"
$searchTerms.Range("D1").Value = "LinkName"
$searchTerms.Range("C1").Value = "SearchType"

# Copy the header formatting (bold font + shaded fill) from an existing sheet.
$wb.Worksheets.Item("TermsPageDrugs").Range("A1:B1").Copy()
$searchTerms.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows.
$searchTerms.Range("A2").Value = "/publications/dictionaries/cancer-drug"
$searchTerms.Range("B2").Value = "Dictionary Page"

$searchTerms.Range("A3").Value = "/publications/dictionaries/cancer-drug?expand=J"
$searchTerms.Range("B3").Value = "Dictionary Page Expand"

$searchTerms.Range("A4").Value = "/publications/dictionaries/cancer-drug/search?contains=false&q=interferon"
$searchTerms.Range("B4").Value = "Dictionary Search Page"

# Column widths (best-fit sized, similar to the sibling dictionary sheets).
$searchTerms.Columns.Item(1).ColumnWidth = 69.45
$searchTerms.Columns.Item(2).ColumnWidth = 21.17
$searchTerms.Columns.Item(3).ColumnWidth = 11.6

[void]$searchTerms.Range("C3").Select()
$searchTerms.Activate()

# ---------------------------------------------------------------------------
# 3. Adjust the selection on TermsPageDrugs to cover the whole data range.
# ---------------------------------------------------------------------------
[void]$wb.Worksheets.Item("TermsPageDrugs").Range("A1:B4").Select()

# Re-activate the SearchTerms sheet so it ends up as the active tab.
$searchTerms.Activate()
